$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.601.48"
$ws.Range("E2").Value = "  +0.93%  "

$ws.Range("D3").Value = "1.896.00"
$ws.Range("E3").Value = "  -0.37%  "

$ws.Range("E4").Value = "  -0.74%  "

$ws.Range("D5").Value = "'247.74"
$ws.Range("E5").Value = "  -2.97%  "

$ws.Range("E6").Value = "  -4.84%  "

$ws.Range("E7").Value = "  -0.83%  "

$ws.Range("D8").Value = "'43.89"
$ws.Range("E8").Value = "  +8.35%  "

$ws.Range("D9").Value = "'0.353"
$ws.Range("E9").Value = "  -4.93%  "

$ws.Range("D10").Value = "'0.0742"
$ws.Range("E10").Value = "  -2.08%  "

$ws.Range("D11").Value = "'0.0971"
$ws.Range("E11").Value = "  -1.70%  "

$ws.Range("D12").Value = "'13.13"
$ws.Range("E12").Value = "  +1.65%  "

$ws.Range("E13").Value = "  -0.35%  "

$ws.Range("D14").Value = "'0.731"
$ws.Range("E14").Value = "  +0.99%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.931.79"
$ws.Range("E15").Value = "  +1.45%  "

$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "'4.96"
$ws.Range("E16").Value = "  +0.17%  "

$ws.Range("D17").Value = "35.557.83"
$ws.Range("E17").Value = "  +0.91%  "

$ws.Range("D18").Value = "'73.83"
$ws.Range("E18").Value = "  -1.12%  "

$ws.Range("E19").Value = "  -2.42%  "

$ws.Range("E20").Value = "  +1.54%  "

$ws.Range("D21").Value = "'12.90"
$ws.Range("E21").Value = "  -0.95%  "

$ws.Range("D22").Value = "'4.96"
$ws.Range("E22").Value = "  -2.60%  "

$ws.Range("E23").Value = "  -0.81%  "

$ws.Range("E24").Value = "  +3.57%  "

$ws.Range("E25").Value = "  -9.27%  "

$ws.Range("D26").Value = "'165.96"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("D27").Value = "'8.53"
$ws.Range("E27").Value = "  -1.53%  "

$ws.Range("D28").Value = "'18.44"
$ws.Range("E28").Value = "  -1.53%  "

$ws.Range("E29").Value = "  -3.65%  "

$ws.Range("D30").Value = "4.128.47"
$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("D31").Value = "'1.77"
$ws.Range("E31").Value = "  +7.95%  "

$ws.Range("D32").Value = "'4.25"
$ws.Range("E32").Value = "  -2.72%  "

$ws.Range("D33").Value = "'0.0583"
$ws.Range("E33").Value = "  -0.86%  "

$ws.Range("D34").Value = "'4.24"
$ws.Range("E34").Value = "  +0.09%  "

$ws.Range("E35").Value = "  -0.81%  "

$ws.Range("D36").Value = "'0.854"
$ws.Range("E36").Value = "  -6.42%  "

$ws.Range("D37").Value = "'2.02"
$ws.Range("E37").Value = "  -0.36%  "

$ws.Range("D38").Value = "'1.54"
$ws.Range("E38").Value = "  -22.50%  "

$ws.Range("E39").Value = "  +5.34%  "

$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "'17.20"
$ws.Range("E40").Value = "  +0.87%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'97.99"
$ws.Range("E41").Value = "  +1.71%  "

$ws.Range("E42").Value = "  -2.13%  "

$ws.Range("E43").Value = "  -2.14%  "

$ws.Range("D44").Value = "1.300.25"
$ws.Range("E44").Value = "  -2.60%  "

$ws.Range("E45").Value = "  -2.32%  "

$ws.Range("D46").Value = "'0.0813"
$ws.Range("E46").Value = "  +7.78%  "

$ws.Range("E47").Value = "  -1.35%  "

$ws.Range("D49").Value = "'12.19"
$ws.Range("E49").Value = "  +4.44%  "

$ws.Range("D50").Value = "'6.36"
$ws.Range("E50").Value = "  -4.88%  "

$ws.Range("D51").Value = "'43.38"
$ws.Range("E51").Value = "  -3.98%  "
